# Changes of 26th May 2022
# Update ShipmentTrackNum (column C) and, where mirrored, PackageTrackNum
# (column D) values for rows 2-22 on Sheet1 with newly generated tracking
# numbers, while keeping the cells stored as text (shared string) with no
# residual style changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param([object]$Worksheet, [string]$CellAddress, [string]$Text)
    $range = $Worksheet.Range($CellAddress)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

$updates = @(
    @{ Row = 2;  Value = "320018638745" },
    @{ Row = 3;  Value = "320018638756" },
    @{ Row = 4;  Value = "320018638789" },
    @{ Row = 5;  Value = "320018638804" },
    @{ Row = 6;  Value = "320018638848" },
    @{ Row = 7;  Value = "320018638860" },
    @{ Row = 8;  Value = "320018638892" },
    @{ Row = 9;  Value = "320018638918" },
    @{ Row = 10; Value = "320018638940" },
    @{ Row = 11; Value = "320018638962" },
    @{ Row = 12; Value = "320018639009" },
    @{ Row = 13; Value = "320018639020" },
    @{ Row = 14; Value = "320018639053" },
    @{ Row = 15; Value = "320018639075" },
    @{ Row = 16; Value = "320018639101" },
    @{ Row = 17; Value = "320018639123" },
    @{ Row = 18; Value = "320018639167" },
    @{ Row = 19; Value = "320018639189" },
    @{ Row = 20; Value = "320018639215" },
    @{ Row = 21; Value = "320018639237" },
    @{ Row = 22; Value = "320018639260" }
)

# Rows where column D previously mirrored column C and must be updated too.
$mirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

foreach ($update in $updates) {
    $row = $update.Row
    $value = $update.Value

    $colCAddress = "C" + $row
    Set-TextValue $ws $colCAddress $value

    if ($mirrorRows -contains $row) {
        $colDAddress = "D" + $row
        Set-TextValue $ws $colDAddress $value
    }
}
